# Update nomenclature. Clean code.
#
# 1. "kite" sheet: remove the "prop.p" row (row 8, value 120) and
#    rename "obgen.p" -> "obGen.p" (row 7, value unchanged).
# 2. "tether" sheet: rename "sigma" -> "sigma_max" (row 9).
# 3. Active sheet/selection bookkeeping: "tether" becomes the active tab
#    (was "gStation"); selections move on "kite" and "tether".

$wb = $excel.ActiveWorkbook

# --- kite sheet -----------------------------------------------------
$kite = $wb.Worksheets.Item("kite")

# Delete the now-obsolete "prop.p" row entirely (row 8: A8="prop.p", B8=120)
$kite.Rows.Item(8).Delete() | Out-Null

# Rename "obgen.p" -> "obGen.p" (still row 7 after the deletion above)
$kite.Range("A7").Value = "obGen.p"

# Move the selection on the kite sheet
$kite.Range("B17").Select() | Out-Null

# --- tether sheet -----------------------------------------------------
$tether = $wb.Worksheets.Item("tether")

# Rename "sigma" -> "sigma_max"
$tether.Range("A9").Value = "sigma_max"

# tether becomes the active/selected sheet (was gStation)
$tether.Activate() | Out-Null
$tether.Range("F9").Select() | Out-Null
